$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 9).Value = 6312
$ws.Cells.Item(3, 9).Value = 6581
$ws.Cells.Item(4, 9).Value = 1509
$ws.Cells.Item(5, 9).Value = 608
$ws.Cells.Item(6, 9).Value = 7455
$ws.Cells.Item(7, 9).Value = 22465
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 9).Value = 176
$ws.Cells.Item(7, 9).Value = 711
$ws.Cells.Item(11, 9).Value = 335
$ws.Cells.Item(14, 9).Value = 128
$ws.Cells.Item(15, 9).Value = 258
$ws.Cells.Item(20, 9).Value = 560
$ws.Cells.Item(25, 9).Value = 121
$ws.Cells.Item(29, 9).Value = 1372
$ws.Cells.Item(31, 9).Value = 226
$ws.Cells.Item(33, 9).Value = 1021
$ws.Cells.Item(36, 9).Value = 306
$ws.Cells.Item(37, 9).Value = 708
$ws.Cells.Item(42, 9).Value = 790
$ws.Cells.Item(45, 9).Value = 47
$ws.Cells.Item(50, 9).Value = 112
$ws.Cells.Item(51, 9).Value = 269
$ws.Cells.Item(52, 9).Value = 481
$ws.Cells.Item(55, 9).Value = 252
$ws.Cells.Item(61, 9).Value = 22
$ws.Cells.Item(63, 9).Value = 78
$ws.Cells.Item(64, 9).Value = 186
$ws.Cells.Item(67, 9).Value = 865
$ws.Cells.Item(73, 9).Value = 206
$ws.Cells.Item(76, 9).Value = 324
$ws.Cells.Item(78, 9).Value = 306
$ws.Cells.Item(79, 9).Value = 638
$ws.Cells.Item(83, 9).Value = 489
$ws.Cells.Item(84, 9).Value = 196
$ws.Cells.Item(85, 9).Value = 1014
$ws.Cells.Item(88, 9).Value = 208
$ws.Cells.Item(89, 9).Value = 265
$ws.Cells.Item(90, 9).Value = 284
$ws.Cells.Item(92, 9).Value = 65
$ws.Cells.Item(93, 9).Value = 128
$ws.Cells.Item(94, 9).Value = 231
$ws.Cells.Item(95, 9).Value = 342
$ws.Cells.Item(96, 9).Value = 246
$ws.Cells.Item(98, 9).Value = 155
$ws.Cells.Item(99, 9).Value = 404
$ws.Cells.Item(101, 9).Value = 22465
$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 9).Value = 289
$ws.Cells.Item(3, 9).Value = 392
$ws.Cells.Item(6, 9).Value = 255
$ws.Cells.Item(7, 9).Value = 1014
$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 9).Value = 129
$ws.Cells.Item(7, 9).Value = 481
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 9).Value = 139
$ws.Cells.Item(7, 9).Value = 335
$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(3, 9).Value = 385
$ws.Cells.Item(4, 9).Value = 85
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(3, 9).Value = 219
$ws.Cells.Item(6, 9).Value = 188
$ws.Cells.Item(7, 9).Value = 711
$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(2, 9).Value = 66
$ws.Cells.Item(6, 9).Value = 91
$ws.Cells.Item(7, 9).Value = 265
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(2, 9).Value = 75
$ws.Cells.Item(7, 9).Value = 246
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(6, 9).Value = 45
$ws.Cells.Item(7, 9).Value = 128
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 9).Value = 212
$ws.Cells.Item(7, 9).Value = 708
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(3, 9).Value = 150
$ws.Cells.Item(7, 9).Value = 404
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 9).Value = 204
$ws.Cells.Item(3, 9).Value = 321
$ws.Cells.Item(7, 9).Value = 865
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(6, 9).Value = 88
$ws.Cells.Item(7, 9).Value = 226
$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(5, 9).Value = 6
$ws.Cells.Item(7, 9).Value = 196
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 9).Value = 166
$ws.Cells.Item(3, 9).Value = 178
$ws.Cells.Item(6, 9).Value = 103
$ws.Cells.Item(7, 9).Value = 489
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(3, 9).Value = 123
$ws.Cells.Item(7, 9).Value = 342
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 9).Value = 231
$ws.Cells.Item(4, 9).Value = 45
$ws.Cells.Item(6, 9).Value = 325
$ws.Cells.Item(7, 9).Value = 1021
$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 9).Value = 473
$ws.Cells.Item(4, 9).Value = 70
$ws.Cells.Item(6, 9).Value = 379
$ws.Cells.Item(7, 9).Value = 1372
$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 9).Value = 210
$ws.Cells.Item(6, 9).Value = 192
$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(3, 9).Value = 75
$ws.Cells.Item(6, 9).Value = 146
$ws.Cells.Item(7, 9).Value = 324
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(3, 9).Value = 245
$ws.Cells.Item(7, 9).Value = 790
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(3, 9).Value = 78
$ws.Cells.Item(5, 9).Value = 7
$ws.Cells.Item(7, 9).Value = 306
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(3, 9).Value = 80
$ws.Cells.Item(7, 9).Value = 252
$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 9).Value = 186
$ws.Cells.Item(3, 9).Value = 205
$ws.Cells.Item(7, 9).Value = 638
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(2, 9).Value = 52
$ws.Cells.Item(4, 9).Value = 13
$ws.Cells.Item(7, 9).Value = 186
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 9).Value = 154
$ws.Cells.Item(3, 9).Value = 159
$ws.Cells.Item(7, 9).Value = 560
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(2, 9).Value = 88
$ws.Cells.Item(6, 9).Value = 96
$ws.Cells.Item(7, 9).Value = 306
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(2, 9).Value = 36
$ws.Cells.Item(7, 9).Value = 128
$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(6, 9).Value = 132
$ws.Cells.Item(7, 9).Value = 231
$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(3, 9).Value = 35
$ws.Cells.Item(7, 9).Value = 121
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 9).Value = 75
$ws.Cells.Item(7, 9).Value = 258
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(3, 9).Value = 13
$ws.Cells.Item(6, 9).Value = 100
$ws.Cells.Item(7, 9).Value = 155
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(2, 9).Value = 33
$ws.Cells.Item(7, 9).Value = 112
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 9).Value = 69
$ws.Cells.Item(7, 9).Value = 206
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(2, 9).Value = 65
$ws.Cells.Item(7, 9).Value = 176
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Cells.Item(2, 9).Value = 21
$ws.Cells.Item(7, 9).Value = 65
$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(6, 9).Value = 63
$ws.Cells.Item(7, 9).Value = 208
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(6, 9).Value = 99
$ws.Cells.Item(7, 9).Value = 284
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(2, 9).Value = 56
$ws.Cells.Item(7, 9).Value = 269
$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(3, 9).Value = 33
$ws.Cells.Item(4, 9).Value = 12
$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Cells.Item(6, 9).Value = 16
$ws.Cells.Item(7, 9).Value = 47
$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Cells.Item(6, 9).Value = 8
$ws.Cells.Item(7, 9).Value = 22
